$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '25.687.51'
Set-TextValue $ws.Range("E2") '  -3.85%  '

Set-TextValue $ws.Range("D3") '1.745.24'
Set-TextValue $ws.Range("E3") '  -5.80%  '

Set-TextValue $ws.Range("D4") '1.001'
Set-TextValue $ws.Range("E4") '  +0.24%  '

Set-TextValue $ws.Range("D5") '236.96'
Set-TextValue $ws.Range("E5") '  -10.11%  '

Set-TextValue $ws.Range("D6") '1.001'
Set-TextValue $ws.Range("E6") '  +0.22%  '

Set-TextValue $ws.Range("D7") '0.4901'
Set-TextValue $ws.Range("E7") '  -8.91%  '

Set-TextValue $ws.Range("D8") '41.56'
Set-TextValue $ws.Range("E8") '  -7.90%  '

Set-TextValue $ws.Range("D9") '0.2501'
Set-TextValue $ws.Range("E9") '  -21.97%  '

Set-TextValue $ws.Range("D10") '0.05937'
Set-TextValue $ws.Range("E10") '  -16.20%  '

Set-TextValue $ws.Range("D11") '1.744.80'
Set-TextValue $ws.Range("E11") '  -5.58%  '

Set-TextValue $ws.Range("D12") '0.06780'
Set-TextValue $ws.Range("E12") '  -13.30%  '

Set-TextValue $ws.Range("D13") '14.74'
Set-TextValue $ws.Range("E13") '  -22.86%  '

Set-TextValue $ws.Range("D14") '4.460'
Set-TextValue $ws.Range("E14") '  -11.84%  '

Set-TextValue $ws.Range("D15") '77.20'
Set-TextValue $ws.Range("E15") '  -14.09%  '

Set-TextValue $ws.Range("D16") '0.5681'
Set-TextValue $ws.Range("E16") '  -27.09%  '

Set-TextValue $ws.Range("D17") '1.001'
Set-TextValue $ws.Range("E17") '  +0.16%  '

Set-TextValue $ws.Range("E18") '  +0.25%  '

Set-TextValue $ws.Range("D19") '25.746.19'
Set-TextValue $ws.Range("E19") '  -3.66%  '

Set-TextValue $ws.Range("D20") '11.48'
Set-TextValue $ws.Range("E20") '  -19.06%  '

Set-TextValue $ws.Range("D21") '0.000006557'
Set-TextValue $ws.Range("E21") '  -18.27%  '

Set-TextValue $ws.Range("D22") '1.962.54'
Set-TextValue $ws.Range("E22") '  -5.92%  '

Set-TextValue $ws.Range("D23") '3.971'
Set-TextValue $ws.Range("E23") '  -14.75%  '

Set-TextValue $ws.Range("D24") '5.027'
Set-TextValue $ws.Range("E24") '  -17.12%  '

Set-TextValue $ws.Range("D25") '7.872'
Set-TextValue $ws.Range("E25") '  -16.76%  '

Set-TextValue $ws.Range("D26") '136.10'
Set-TextValue $ws.Range("E26") '  -4.67%  '

Set-TextValue $ws.Range("D27") '1.475'
Set-TextValue $ws.Range("E27") '  -13.67%  '

Set-TextValue $ws.Range("D28") '1.823'
Set-TextValue $ws.Range("E28") '  -18.35%  '

Set-TextValue $ws.Range("D29") '14.63'
Set-TextValue $ws.Range("E29") '  -14.78%  '

Set-TextValue $ws.Range("D30") '101.94'
Set-TextValue $ws.Range("E30") '  -9.01%  '

Set-TextValue $ws.Range("D31") '3.767'
Set-TextValue $ws.Range("E31") '  -12.54%  '

Set-TextValue $ws.Range("D32") '0.08072'
Set-TextValue $ws.Range("E32") '  -7.97%  '

Set-TextValue $ws.Range("D33") '3.324'
Set-TextValue $ws.Range("E33") '  -19.65%  '

Set-TextValue $ws.Range("D34") '0.04394'
Set-TextValue $ws.Range("E34") '  -10.11%  '

Set-TextValue $ws.Range("E35") '  +0.23%  '

Set-TextValue $ws.Range("D36") '2.612'
Set-TextValue $ws.Range("E36") '  -8.92%  '

Set-TextValue $ws.Range("D37") '0.9737'
Set-TextValue $ws.Range("E37") '  -15.22%  '

Set-TextValue $ws.Range("D38") '0.5998'
Set-TextValue $ws.Range("E38") '  -18.80%  '

Set-TextValue $ws.Range("D39") '2.681'
Set-TextValue $ws.Range("E39") '  -13.84%  '

Set-TextValue $ws.Range("D40") '2.025'
Set-TextValue $ws.Range("E40") '  -14.71%  '

Set-TextValue $ws.Range("D41") '1.001'
Set-TextValue $ws.Range("E41") '  +0.22%  '

Set-TextValue $ws.Range("D42") '103.61'
Set-TextValue $ws.Range("E42") '  -5.43%  '

Set-TextValue $ws.Range("D43") '0.01494'
Set-TextValue $ws.Range("E43") '  -14.90%  '

Set-TextValue $ws.Range("D44") '0.7590'
Set-TextValue $ws.Range("E44") '  -16.93%  '

Set-TextValue $ws.Range("D45") '5.162'
Set-TextValue $ws.Range("E45") '  -13.05%  '

Set-TextValue $ws.Range("D46") '0.3715'
Set-TextValue $ws.Range("E46") '  -23.64%  '

Set-TextValue $ws.Range("D47") '0.05111'
Set-TextValue $ws.Range("E47") '  -12.59%  '

Set-TextValue $ws.Range("D48") '0.1071'
Set-TextValue $ws.Range("E48") '  -15.09%  '

Set-TextValue $ws.Range("D49") '30.14'
Set-TextValue $ws.Range("E49") '  -14.26%  '

Set-TextValue $ws.Range("D50") '52.52'
Set-TextValue $ws.Range("E50") '  -13.21%  '

Set-TextValue $ws.Range("D51") '5.868'
Set-TextValue $ws.Range("E51") '  -24.37%  '
